$wb = $excel.ActiveWorkbook

# --- Worklist sheet (sheet1): update Dest (C) and DestWell (D) columns ---
$wsWorklist = $wb.Worksheets.Item("Worklist")

$wsWorklist.Range("C18").Value = "384-2"
$wsWorklist.Range("D18").Value = "33"
$wsWorklist.Range("C19").Value = "384-2"
$wsWorklist.Range("D19").Value = "35"
$wsWorklist.Range("C20").Value = "384-2"
$wsWorklist.Range("D20").Value = "37"
$wsWorklist.Range("C21").Value = "384-2"
$wsWorklist.Range("D21").Value = "39"
$wsWorklist.Range("C22").Value = "384-2"
$wsWorklist.Range("D22").Value = "41"
$wsWorklist.Range("C23").Value = "384-2"
$wsWorklist.Range("D23").Value = "43"
$wsWorklist.Range("C24").Value = "384-2"
$wsWorklist.Range("D24").Value = "45"
$wsWorklist.Range("C25").Value = "384-2"
$wsWorklist.Range("D25").Value = "47"
$wsWorklist.Range("C26").Value = "384-2"
$wsWorklist.Range("D26").Value = "34"
$wsWorklist.Range("C27").Value = "384-2"
$wsWorklist.Range("D27").Value = "36"
$wsWorklist.Range("C28").Value = "384-2"
$wsWorklist.Range("D28").Value = "38"
$wsWorklist.Range("C29").Value = "384-2"
$wsWorklist.Range("D29").Value = "40"
$wsWorklist.Range("C30").Value = "384-2"
$wsWorklist.Range("D30").Value = "42"
$wsWorklist.Range("C31").Value = "384-2"
$wsWorklist.Range("D31").Value = "44"
$wsWorklist.Range("C32").Value = "384-2"
$wsWorklist.Range("D32").Value = "46"
$wsWorklist.Range("C33").Value = "384-2"
$wsWorklist.Range("D33").Value = "48"
$wsWorklist.Range("C34").Value = "384-3"
$wsWorklist.Range("C35").Value = "384-3"
$wsWorklist.Range("C36").Value = "384-3"
$wsWorklist.Range("C37").Value = "384-3"
$wsWorklist.Range("C38").Value = "384-3"
$wsWorklist.Range("C39").Value = "384-3"
$wsWorklist.Range("C40").Value = "384-3"
$wsWorklist.Range("C41").Value = "384-3"
$wsWorklist.Range("C42").Value = "384-3"
$wsWorklist.Range("C43").Value = "384-3"
$wsWorklist.Range("C44").Value = "384-3"
$wsWorklist.Range("C45").Value = "384-3"
$wsWorklist.Range("C46").Value = "384-3"
$wsWorklist.Range("C47").Value = "384-3"
$wsWorklist.Range("C48").Value = "384-3"
$wsWorklist.Range("C49").Value = "384-3"
$wsWorklist.Range("C50").Value = "384-4"
$wsWorklist.Range("D50").Value = "33"
$wsWorklist.Range("C51").Value = "384-4"
$wsWorklist.Range("D51").Value = "35"
$wsWorklist.Range("C52").Value = "384-4"
$wsWorklist.Range("D52").Value = "37"
$wsWorklist.Range("C53").Value = "384-4"
$wsWorklist.Range("D53").Value = "39"
$wsWorklist.Range("C54").Value = "384-4"
$wsWorklist.Range("D54").Value = "41"
$wsWorklist.Range("C55").Value = "384-4"
$wsWorklist.Range("D55").Value = "43"
$wsWorklist.Range("C56").Value = "384-4"
$wsWorklist.Range("D56").Value = "45"
$wsWorklist.Range("C57").Value = "384-4"
$wsWorklist.Range("D57").Value = "47"
$wsWorklist.Range("C58").Value = "384-4"
$wsWorklist.Range("D58").Value = "34"
$wsWorklist.Range("C59").Value = "384-4"
$wsWorklist.Range("D59").Value = "36"
$wsWorklist.Range("C60").Value = "384-4"
$wsWorklist.Range("D60").Value = "38"
$wsWorklist.Range("C61").Value = "384-4"
$wsWorklist.Range("D61").Value = "40"
$wsWorklist.Range("C62").Value = "384-4"
$wsWorklist.Range("D62").Value = "42"
$wsWorklist.Range("C63").Value = "384-4"
$wsWorklist.Range("D63").Value = "44"
$wsWorklist.Range("C64").Value = "384-4"
$wsWorklist.Range("D64").Value = "46"
$wsWorklist.Range("C65").Value = "384-4"
$wsWorklist.Range("D65").Value = "48"
$wsWorklist.Range("C66").Value = "384-5"
$wsWorklist.Range("C67").Value = "384-5"
$wsWorklist.Range("C68").Value = "384-5"
$wsWorklist.Range("C69").Value = "384-5"
$wsWorklist.Range("C70").Value = "384-5"
$wsWorklist.Range("C71").Value = "384-5"
$wsWorklist.Range("C72").Value = "384-5"
$wsWorklist.Range("C73").Value = "384-5"
$wsWorklist.Range("C74").Value = "384-5"
$wsWorklist.Range("C75").Value = "384-5"
$wsWorklist.Range("C76").Value = "384-5"
$wsWorklist.Range("C77").Value = "384-5"
$wsWorklist.Range("C78").Value = "384-5"
$wsWorklist.Range("C79").Value = "384-5"
$wsWorklist.Range("C80").Value = "384-5"
$wsWorklist.Range("C81").Value = "384-5"
$wsWorklist.Range("C82").Value = "384-6"
$wsWorklist.Range("D82").Value = "33"
$wsWorklist.Range("C83").Value = "384-6"
$wsWorklist.Range("D83").Value = "35"
$wsWorklist.Range("C84").Value = "384-6"
$wsWorklist.Range("D84").Value = "37"
$wsWorklist.Range("C85").Value = "384-6"
$wsWorklist.Range("D85").Value = "39"
$wsWorklist.Range("C86").Value = "384-6"
$wsWorklist.Range("D86").Value = "41"
$wsWorklist.Range("C87").Value = "384-6"
$wsWorklist.Range("D87").Value = "43"
$wsWorklist.Range("C88").Value = "384-6"
$wsWorklist.Range("D88").Value = "45"
$wsWorklist.Range("C89").Value = "384-6"
$wsWorklist.Range("D89").Value = "47"
$wsWorklist.Range("C90").Value = "384-6"
$wsWorklist.Range("D90").Value = "34"
$wsWorklist.Range("C91").Value = "384-6"
$wsWorklist.Range("D91").Value = "36"
$wsWorklist.Range("C92").Value = "384-6"
$wsWorklist.Range("D92").Value = "38"
$wsWorklist.Range("C93").Value = "384-6"
$wsWorklist.Range("D93").Value = "40"
$wsWorklist.Range("C94").Value = "384-6"
$wsWorklist.Range("D94").Value = "42"
$wsWorklist.Range("C95").Value = "384-6"
$wsWorklist.Range("D95").Value = "44"
$wsWorklist.Range("C96").Value = "384-6"
$wsWorklist.Range("D96").Value = "46"
$wsWorklist.Range("C97").Value = "384-6"
$wsWorklist.Range("D97").Value = "48"
$wsWorklist.Range("C98").Value = "384-7"
$wsWorklist.Range("C99").Value = "384-7"
$wsWorklist.Range("C100").Value = "384-7"
$wsWorklist.Range("C101").Value = "384-7"
$wsWorklist.Range("C102").Value = "384-7"
$wsWorklist.Range("C103").Value = "384-7"
$wsWorklist.Range("C104").Value = "384-7"
$wsWorklist.Range("C105").Value = "384-7"
$wsWorklist.Range("C106").Value = "384-7"
$wsWorklist.Range("C107").Value = "384-7"
$wsWorklist.Range("C108").Value = "384-7"
$wsWorklist.Range("C109").Value = "384-7"
$wsWorklist.Range("C110").Value = "384-7"
$wsWorklist.Range("C111").Value = "384-7"
$wsWorklist.Range("C112").Value = "384-7"
$wsWorklist.Range("C113").Value = "384-7"
$wsWorklist.Range("C114").Value = "384-8"
$wsWorklist.Range("D114").Value = "33"
$wsWorklist.Range("C115").Value = "384-8"
$wsWorklist.Range("D115").Value = "34"
$wsWorklist.Range("C116").Value = "384-8"
$wsWorklist.Range("D116").Value = "35"

# --- Platemap sheet (sheet2): update Dest (D), DestWell (E), Well (F) columns ---
$wsPlatemap = $wb.Worksheets.Item("Platemap")

$wsPlatemap.Range("D17").Value = "384-2"
$wsPlatemap.Range("E17").Value = "33"
$wsPlatemap.Range("F17").Value = "A03"
$wsPlatemap.Range("D18").Value = "384-2"
$wsPlatemap.Range("E18").Value = "35"
$wsPlatemap.Range("F18").Value = "C03"
$wsPlatemap.Range("D19").Value = "384-2"
$wsPlatemap.Range("E19").Value = "37"
$wsPlatemap.Range("F19").Value = "E03"
$wsPlatemap.Range("D20").Value = "384-2"
$wsPlatemap.Range("E20").Value = "39"
$wsPlatemap.Range("F20").Value = "G03"
$wsPlatemap.Range("D21").Value = "384-2"
$wsPlatemap.Range("E21").Value = "41"
$wsPlatemap.Range("F21").Value = "I03"
$wsPlatemap.Range("D22").Value = "384-2"
$wsPlatemap.Range("E22").Value = "43"
$wsPlatemap.Range("F22").Value = "K03"
$wsPlatemap.Range("D23").Value = "384-2"
$wsPlatemap.Range("E23").Value = "45"
$wsPlatemap.Range("F23").Value = "M03"
$wsPlatemap.Range("D24").Value = "384-2"
$wsPlatemap.Range("E24").Value = "47"
$wsPlatemap.Range("F24").Value = "O03"
$wsPlatemap.Range("D25").Value = "384-2"
$wsPlatemap.Range("E25").Value = "34"
$wsPlatemap.Range("F25").Value = "B03"
$wsPlatemap.Range("D26").Value = "384-2"
$wsPlatemap.Range("E26").Value = "36"
$wsPlatemap.Range("F26").Value = "D03"
$wsPlatemap.Range("D27").Value = "384-2"
$wsPlatemap.Range("E27").Value = "38"
$wsPlatemap.Range("F27").Value = "F03"
$wsPlatemap.Range("D28").Value = "384-2"
$wsPlatemap.Range("E28").Value = "40"
$wsPlatemap.Range("F28").Value = "H03"
$wsPlatemap.Range("D29").Value = "384-2"
$wsPlatemap.Range("E29").Value = "42"
$wsPlatemap.Range("F29").Value = "J03"
$wsPlatemap.Range("D30").Value = "384-2"
$wsPlatemap.Range("E30").Value = "44"
$wsPlatemap.Range("F30").Value = "L03"
$wsPlatemap.Range("D31").Value = "384-2"
$wsPlatemap.Range("E31").Value = "46"
$wsPlatemap.Range("F31").Value = "N03"
$wsPlatemap.Range("D32").Value = "384-2"
$wsPlatemap.Range("E32").Value = "48"
$wsPlatemap.Range("F32").Value = "P03"
$wsPlatemap.Range("D33").Value = "384-3"
$wsPlatemap.Range("D34").Value = "384-3"
$wsPlatemap.Range("D35").Value = "384-3"
$wsPlatemap.Range("D36").Value = "384-3"
$wsPlatemap.Range("D37").Value = "384-3"
$wsPlatemap.Range("D38").Value = "384-3"
$wsPlatemap.Range("D39").Value = "384-3"
$wsPlatemap.Range("D40").Value = "384-3"
$wsPlatemap.Range("D41").Value = "384-3"
$wsPlatemap.Range("D42").Value = "384-3"
$wsPlatemap.Range("D43").Value = "384-3"
$wsPlatemap.Range("D44").Value = "384-3"
$wsPlatemap.Range("D45").Value = "384-3"
$wsPlatemap.Range("D46").Value = "384-3"
$wsPlatemap.Range("D47").Value = "384-3"
$wsPlatemap.Range("D48").Value = "384-3"
$wsPlatemap.Range("D49").Value = "384-4"
$wsPlatemap.Range("E49").Value = "33"
$wsPlatemap.Range("F49").Value = "A03"
$wsPlatemap.Range("D50").Value = "384-4"
$wsPlatemap.Range("E50").Value = "35"
$wsPlatemap.Range("F50").Value = "C03"
$wsPlatemap.Range("D51").Value = "384-4"
$wsPlatemap.Range("E51").Value = "37"
$wsPlatemap.Range("F51").Value = "E03"
$wsPlatemap.Range("D52").Value = "384-4"
$wsPlatemap.Range("E52").Value = "39"
$wsPlatemap.Range("F52").Value = "G03"
$wsPlatemap.Range("D53").Value = "384-4"
$wsPlatemap.Range("E53").Value = "41"
$wsPlatemap.Range("F53").Value = "I03"
$wsPlatemap.Range("D54").Value = "384-4"
$wsPlatemap.Range("E54").Value = "43"
$wsPlatemap.Range("F54").Value = "K03"
$wsPlatemap.Range("D55").Value = "384-4"
$wsPlatemap.Range("E55").Value = "45"
$wsPlatemap.Range("F55").Value = "M03"
$wsPlatemap.Range("D56").Value = "384-4"
$wsPlatemap.Range("E56").Value = "47"
$wsPlatemap.Range("F56").Value = "O03"
$wsPlatemap.Range("D57").Value = "384-4"
$wsPlatemap.Range("E57").Value = "34"
$wsPlatemap.Range("F57").Value = "B03"
$wsPlatemap.Range("D58").Value = "384-4"
$wsPlatemap.Range("E58").Value = "36"
$wsPlatemap.Range("F58").Value = "D03"
$wsPlatemap.Range("D59").Value = "384-4"
$wsPlatemap.Range("E59").Value = "38"
$wsPlatemap.Range("F59").Value = "F03"
$wsPlatemap.Range("D60").Value = "384-4"
$wsPlatemap.Range("E60").Value = "40"
$wsPlatemap.Range("F60").Value = "H03"
$wsPlatemap.Range("D61").Value = "384-4"
$wsPlatemap.Range("E61").Value = "42"
$wsPlatemap.Range("F61").Value = "J03"
$wsPlatemap.Range("D62").Value = "384-4"
$wsPlatemap.Range("E62").Value = "44"
$wsPlatemap.Range("F62").Value = "L03"
$wsPlatemap.Range("D63").Value = "384-4"
$wsPlatemap.Range("E63").Value = "46"
$wsPlatemap.Range("F63").Value = "N03"
$wsPlatemap.Range("D64").Value = "384-4"
$wsPlatemap.Range("E64").Value = "48"
$wsPlatemap.Range("F64").Value = "P03"
$wsPlatemap.Range("D65").Value = "384-5"
$wsPlatemap.Range("D66").Value = "384-5"
$wsPlatemap.Range("D67").Value = "384-5"
$wsPlatemap.Range("D68").Value = "384-5"
$wsPlatemap.Range("D69").Value = "384-5"
$wsPlatemap.Range("D70").Value = "384-5"
$wsPlatemap.Range("D71").Value = "384-5"
$wsPlatemap.Range("D72").Value = "384-5"
$wsPlatemap.Range("D73").Value = "384-5"
$wsPlatemap.Range("D74").Value = "384-5"
$wsPlatemap.Range("D75").Value = "384-5"
$wsPlatemap.Range("D76").Value = "384-5"
$wsPlatemap.Range("D77").Value = "384-5"
$wsPlatemap.Range("D78").Value = "384-5"
$wsPlatemap.Range("D79").Value = "384-5"
$wsPlatemap.Range("D80").Value = "384-5"
$wsPlatemap.Range("D81").Value = "384-6"
$wsPlatemap.Range("E81").Value = "33"
$wsPlatemap.Range("F81").Value = "A03"
$wsPlatemap.Range("D82").Value = "384-6"
$wsPlatemap.Range("E82").Value = "35"
$wsPlatemap.Range("F82").Value = "C03"
$wsPlatemap.Range("D83").Value = "384-6"
$wsPlatemap.Range("E83").Value = "37"
$wsPlatemap.Range("F83").Value = "E03"
$wsPlatemap.Range("D84").Value = "384-6"
$wsPlatemap.Range("E84").Value = "39"
$wsPlatemap.Range("F84").Value = "G03"
$wsPlatemap.Range("D85").Value = "384-6"
$wsPlatemap.Range("E85").Value = "41"
$wsPlatemap.Range("F85").Value = "I03"
$wsPlatemap.Range("D86").Value = "384-6"
$wsPlatemap.Range("E86").Value = "43"
$wsPlatemap.Range("F86").Value = "K03"
$wsPlatemap.Range("D87").Value = "384-6"
$wsPlatemap.Range("E87").Value = "45"
$wsPlatemap.Range("F87").Value = "M03"
$wsPlatemap.Range("D88").Value = "384-6"
$wsPlatemap.Range("E88").Value = "47"
$wsPlatemap.Range("F88").Value = "O03"
$wsPlatemap.Range("D89").Value = "384-6"
$wsPlatemap.Range("E89").Value = "34"
$wsPlatemap.Range("F89").Value = "B03"
$wsPlatemap.Range("D90").Value = "384-6"
$wsPlatemap.Range("E90").Value = "36"
$wsPlatemap.Range("F90").Value = "D03"
$wsPlatemap.Range("D91").Value = "384-6"
$wsPlatemap.Range("E91").Value = "38"
$wsPlatemap.Range("F91").Value = "F03"
$wsPlatemap.Range("D92").Value = "384-6"
$wsPlatemap.Range("E92").Value = "40"
$wsPlatemap.Range("F92").Value = "H03"
$wsPlatemap.Range("D93").Value = "384-6"
$wsPlatemap.Range("E93").Value = "42"
$wsPlatemap.Range("F93").Value = "J03"
$wsPlatemap.Range("D94").Value = "384-6"
$wsPlatemap.Range("E94").Value = "44"
$wsPlatemap.Range("F94").Value = "L03"
$wsPlatemap.Range("D95").Value = "384-6"
$wsPlatemap.Range("E95").Value = "46"
$wsPlatemap.Range("F95").Value = "N03"
$wsPlatemap.Range("D96").Value = "384-6"
$wsPlatemap.Range("E96").Value = "48"
$wsPlatemap.Range("F96").Value = "P03"
$wsPlatemap.Range("D97").Value = "384-7"
$wsPlatemap.Range("D98").Value = "384-7"
$wsPlatemap.Range("D99").Value = "384-7"
$wsPlatemap.Range("D100").Value = "384-7"
$wsPlatemap.Range("D101").Value = "384-7"
$wsPlatemap.Range("D102").Value = "384-7"
$wsPlatemap.Range("D103").Value = "384-7"
$wsPlatemap.Range("D104").Value = "384-7"
$wsPlatemap.Range("D105").Value = "384-7"
$wsPlatemap.Range("D106").Value = "384-7"
$wsPlatemap.Range("D107").Value = "384-7"
$wsPlatemap.Range("D108").Value = "384-7"
$wsPlatemap.Range("D109").Value = "384-7"
$wsPlatemap.Range("D110").Value = "384-7"
$wsPlatemap.Range("D111").Value = "384-7"
$wsPlatemap.Range("D112").Value = "384-7"
$wsPlatemap.Range("D113").Value = "384-8"
$wsPlatemap.Range("E113").Value = "33"
$wsPlatemap.Range("F113").Value = "A03"
$wsPlatemap.Range("D114").Value = "384-8"
$wsPlatemap.Range("E114").Value = "34"
$wsPlatemap.Range("F114").Value = "B03"
$wsPlatemap.Range("D115").Value = "384-8"
$wsPlatemap.Range("E115").Value = "35"
$wsPlatemap.Range("F115").Value = "C03"
